$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.9693284034729
$ws.Range("B1").Value = 6.331477642059326
$ws.Range("C1").Value = 3.323112487792969
$ws.Range("D1").Value = 1.465604424476624
$ws.Range("E1").Value = 1.030220150947571
